$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current last row (140), pushing the
# existing last row (Red Beaut / Primera) down to row 142.
$ws.Range("A140:T141").Insert()

# New row 140: Lemon / Primera
$ws.Range("A140").Value = 4
$ws.Range("B140").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C140").Value = "Los Lagos"
$ws.Range("D140").Value = 44595
$ws.Range("D140").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E140").Value = 10
$ws.Range("F140").Value = "Fruta"
$ws.Range("G140").Value = 100103
$ws.Range("H140").Value = "Frutos de hueso (carozo)"
$ws.Range("I140").Value = 100103002
$ws.Range("J140").Value = "Ciruela"
$ws.Range("K140").Value = "Lemon"
$ws.Range("L140").Value = "Primera"
$ws.Range("M140").Value = 200
$ws.Range("N140").Value = 15000
$ws.Range("O140").Value = 16000
$ws.Range("P140").Value = 15500
$ws.Range("Q140").Value = "`$/caja 15 kilos granel"
$ws.Range("R140").Value = "Región de O'Higgins"
$ws.Range("S140").Value = 1033
$ws.Range("T140").Value = 15

# New row 141: Lemon / Segunda
$ws.Range("A141").Value = 4
$ws.Range("B141").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C141").Value = "Los Lagos"
$ws.Range("D141").Value = 44595
$ws.Range("D141").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E141").Value = 10
$ws.Range("F141").Value = "Fruta"
$ws.Range("G141").Value = 100103
$ws.Range("H141").Value = "Frutos de hueso (carozo)"
$ws.Range("I141").Value = 100103002
$ws.Range("J141").Value = "Ciruela"
$ws.Range("K141").Value = "Lemon"
$ws.Range("L141").Value = "Segunda"
$ws.Range("M141").Value = 150
$ws.Range("N141").Value = 13000
$ws.Range("O141").Value = 13000
$ws.Range("P141").Value = 13000
$ws.Range("Q141").Value = "`$/caja 15 kilos granel"
$ws.Range("R141").Value = "Región de O'Higgins"
$ws.Range("S141").Value = 867
$ws.Range("T141").Value = 15

"done"
